# Add "2022-Q4" quarterly data to the workbook:
#  1. Insert a new detail worksheet "2022-Q4" right after "总计", populated with
#     the fund-level holdings for that quarter.
#  2. Insert a corresponding summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet positioned right after "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Use "2022-Q2" as a formatting template (same layout as every other quarter sheet)
$template = $wb.Worksheets.Item("2022-Q2")

# Header row formatting (columns B:H carry the bold/bordered header style)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Column A formatting for the two data rows (bold/bordered index style)
$template.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

# Header labels
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B:G hold text values (fund codes / numbers-as-text), force text storage
$newSheet.Range("B2:G3").NumberFormat = "@"

# Row 2 - 华安沪港深外延增长混合A
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "001694"
$newSheet.Range("C2").Value = "华安沪港深外延增长混合A"
$newSheet.Range("D2").Value = "38.21"
$newSheet.Range("E2").Value = "94.15"
$newSheet.Range("F2").Value = "2.24"
$newSheet.Range("G2").Value = "0.8559"
$newSheet.Range("H2").Value = 10

# Row 3 - 华安沪港深外延增长混合C
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "014972"
$newSheet.Range("C3").Value = "华安沪港深外延增长混合C"
$newSheet.Range("D3").Value = "0.13"
$newSheet.Range("E3").Value = "94.15"
$newSheet.Range("F3").Value = "2.24"
$newSheet.Range("G3").Value = "0.0029"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2. Insert the matching summary row into the "总计" sheet (row 2, right under
#    the header), pushing all existing quarters down by one row.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Re-apply the index style (column A) and carry over the plain data-cell
# formatting from the row that was just pushed down to row 3.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)    # xlPasteFormats
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.86

# Renumber the index column (A) for the remaining (shifted) quarter rows
for ($r = 3; $r -le 8; $r++) {
    $totalSheet.Range("A$r").Value = $r - 2
}

Write-Host "2022-Q4 data added"
